$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.649.08"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.292.71"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "2.636.82"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "2.294.90"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "43.668.44"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.242"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "2.516.18"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.428"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
